# Se procesan de nuevo los datos con las nuevas dimensiones curadas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (personas-mano-obra-familiar-conyuges-...-secundaria) was curated
# from a dimension to a measure.
$ws.Range("A2").Value = "iaest-measure:personas-mano-obra-familiar-conyuges-que-trabajan-en-otra-actividad-como-secundaria"
$ws.Range("A3").Value = "medida"
$ws.Range("A4").Value = "xsd:int"
$ws.Range("A5").Clear()

# Column F (provincia) was curated from the generic sdmx refArea dimension
# to a plain measure.
$ws.Range("F2").Value = "iaest-measure:provincia"
$ws.Range("F3").Value = "medida"
$ws.Range("F4").Value = "xsd:int"
